{"js": "// Update the two-digit multiplication problems in the table to match\n// the target revision. Each cell's \"A\u00d7B=\" text is replaced with its new\n// value; old values are unique in the document, so a single search per\n// pair is unambiguous.\nconst pairs = [\n  [\"71\u00d776=\", \"71\u00d770=\"],\n  [\"67\u00d715=\", \"37\u00d798=\"],\n  [\"74\u00d775=\", \"30\u00d779=\"],\n  [\"94\u00d741=\", \"12\u00d717=\"],\n  [\"42\u00d773=\", \"75\u00d743=\"],\n  [\"76\u00d787=\", \"95\u00d791=\"],\n  [\"78\u00d779=\", \"72\u00d747=\"],\n  [\"68\u00d778=\", \"49\u00d778=\"],\n  [\"40\u00d766=\", \"81\u00d745=\"],\n  [\"32\u00d776=\", \"99\u00d781=\"],\n  [\"39\u00d777=\", \"71\u00d785=\"],\n  [\"98\u00d711=\", \"56\u00d712=\"],\n  [\"62\u00d714=\", \"20\u00d713=\"],\n  [\"22\u00d766=\", \"29\u00d783=\"],\n  [\"14\u00d745=\", \"84\u00d794=\"],\n  [\"70\u00d744=\", \"59\u00d731=\"],\n  [\"61\u00d797=\", \"21\u00d734=\"],\n  [\"19\u00d761=\", \"50\u00d793=\"],\n  [\"78\u00d764=\", \"21\u00d723=\"],\n  [\"15\u00d728=\", \"14\u00d740=\"],\n  [\"94\u00d723=\", \"67\u00d785=\"],\n  [\"48\u00d714=\", \"74\u00d752=\"],\n  [\"36\u00d710=\", \"63\u00d772=\"],\n  [\"60\u00d755=\", \"45\u00d714=\"],\n  [\"38\u00d788=\", \"52\u00d720=\"],\n  [\"56\u00d728=\", \"52\u00d744=\"],\n  [\"80\u00d729=\", \"58\u00d733=\"],\n  [\"95\u00d765=\", \"88\u00d786=\"],\n  [\"86\u00d721=\", \"52\u00d722=\"],\n  [\"48\u00d796=\", \"39\u00d789=\"],\n  [\"34\u00d753=\", \"15\u00d766=\"],\n  [\"44\u00d783=\", \"68\u00d788=\"],\n  [\"87\u00d786=\", \"21\u00d797=\"],\n  [\"21\u00d789=\", \"99\u00d732=\"],\n  [\"40\u00d777=\", \"42\u00d763=\"],\n  [\"98\u00d762=\", \"97\u00d766=\"],\n  [\"36\u00d760=\", \"96\u00d779=\"],\n  [\"90\u00d745=\", \"13\u00d783=\"],\n  [\"83\u00d773=\", \"72\u00d760=\"],\n  [\"87\u00d724=\", \"50\u00d771=\"],\n  [\"42\u00d795=\", \"17\u00d710=\"],\n  [\"61\u00d736=\", \"61\u00d787=\"],\n  [\"74\u00d745=\", \"43\u00d748=\"],\n  [\"100\u00d730=\", \"80\u00d723=\"],\n  [\"31\u00d726=\", \"22\u00d757=\"],\n  [\"24\u00d770=\", \"18\u00d733=\"],\n  [\"74\u00d751=\", \"80\u00d789=\"],\n  [\"47\u00d744=\", \"13\u00d798=\"],\n  [\"46\u00d727=\", \"10\u00d795=\"],\n  [\"68\u00d727=\", \"22\u00d714=\"],\n  [\"27\u00d770=\", \"91\u00d779=\"],\n  [\"72\u00d777=\", \"85\u00d738=\"],\n  [\"53\u00d732=\", \"95\u00d710=\"],\n  [\"94\u00d770=\", \"11\u00d769=\"],\n  [\"21\u00d748=\", \"41\u00d719=\"],\n  [\"87\u00d743=\", \"79\u00d754=\"],\n  [\"32\u00d783=\", \"86\u00d781=\"],\n  [\"95\u00d714=\", \"94\u00d794=\"],\n  [\"71\u00d738=\", \"69\u00d793=\"],\n  [\"79\u00d798=\", \"30\u00d789=\"],\n  [\"69\u00d741=\", \"76\u00d779=\"],\n  [\"24\u00d740=\", \"74\u00d765=\"],\n  [\"30\u00d761=\", \"61\u00d787=\"],\n  [\"30\u00d773=\", \"13\u00d777=\"],\n  [\"13\u00d769=\", \"24\u00d721=\"],\n  [\"15\u00d721=\", \"73\u00d757=\"],\n  [\"32\u00d768=\", \"74\u00d730=\"],\n  [\"90\u00d762=\", \"24\u00d735=\"],\n  [\"62\u00d747=\", \"65\u00d745=\"],\n  [\"54\u00d726=\", \"30\u00d796=\"],\n  [\"42\u00d750=\", \"72\u00d780=\"],\n  [\"42\u00d729=\", \"63\u00d740=\"],\n  [\"12\u00d734=\", \"28\u00d762=\"],\n  [\"22\u00d734=\", \"38\u00d799=\"],\n  [\"41\u00d760=\", \"93\u00d718=\"],\n  [\"92\u00d782=\", \"81\u00d738=\"],\n  [\"76\u00d776=\", \"47\u00d726=\"],\n  [\"70\u00d732=\", \"19\u00d770=\"],\n  [\"16\u00d753=\", \"59\u00d787=\"],\n  [\"98\u00d757=\", \"95\u00d790=\"],\n  [\"74\u00d738=\", \"42\u00d718=\"],\n  [\"56\u00d725=\", \"30\u00d731=\"],\n  [\"54\u00d7100=\", \"23\u00d765=\"],\n  [\"52\u00d727=\", \"55\u00d764=\"],\n  [\"44\u00d722=\", \"72\u00d786=\"],\n  [\"87\u00d733=\", \"81\u00d726=\"],\n  [\"11\u00d713=\", \"59\u00d758=\"],\n  [\"89\u00d778=\", \"81\u00d739=\"],\n  [\"38\u00d782=\", \"53\u00d748=\"],\n  [\"50\u00d760=\", \"66\u00d756=\"],\n  [\"42\u00d783=\", \"23\u00d763=\"],\n  [\"11\u00d745=\", \"97\u00d787=\"],\n  [\"67\u00d769=\", \"24\u00d743=\"],\n  [\"12\u00d757=\", \"66\u00d721=\"],\n  [\"22\u00d772=\", \"40\u00d790=\"],\n  [\"44\u00d792=\", \"21\u00d733=\"],\n  [\"11\u00d743=\", \"63\u00d751=\"],\n  [\"55\u00d722=\", \"52\u00d767=\"],\n  [\"83\u00d768=\", \"21\u00d736=\"],\n  [\"30\u00d712=\", \"24\u00d714=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldText}\", found ${results.items.length}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update multiplication problems in the table per the commit diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('71\u00d776=', '71\u00d770='),\n    @('67\u00d715=', '37\u00d798='),\n    @('74\u00d775=', '30\u00d779='),\n    @('94\u00d741=', '12\u00d717='),\n    @('42\u00d773=', '75\u00d743='),\n    @('76\u00d787=', '95\u00d791='),\n    @('78\u00d779=', '72\u00d747='),\n    @('68\u00d778=', '49\u00d778='),\n    @('40\u00d766=', '81\u00d745='),\n    @('32\u00d776=', '99\u00d781='),\n    @('39\u00d777=', '71\u00d785='),\n    @('98\u00d711=', '56\u00d712='),\n    @('62\u00d714=', '20\u00d713='),\n    @('22\u00d766=', '29\u00d783='),\n    @('14\u00d745=', '84\u00d794='),\n    @('70\u00d744=', '59\u00d731='),\n    @('61\u00d797=', '21\u00d734='),\n    @('19\u00d761=', '50\u00d793='),\n    @('78\u00d764=', '21\u00d723='),\n    @('15\u00d728=', '14\u00d740='),\n    @('94\u00d723=', '67\u00d785='),\n    @('48\u00d714=', '74\u00d752='),\n    @('36\u00d710=', '63\u00d772='),\n    @('60\u00d755=', '45\u00d714='),\n    @('38\u00d788=', '52\u00d720='),\n    @('56\u00d728=', '52\u00d744='),\n    @('80\u00d729=', '58\u00d733='),\n    @('95\u00d765=', '88\u00d786='),\n    @('86\u00d721=', '52\u00d722='),\n    @('48\u00d796=', '39\u00d789='),\n    @('34\u00d753=', '15\u00d766='),\n    @('44\u00d783=', '68\u00d788='),\n    @('87\u00d786=', '21\u00d797='),\n    @('21\u00d789=', '99\u00d732='),\n    @('40\u00d777=', '42\u00d763='),\n    @('98\u00d762=', '97\u00d766='),\n    @('36\u00d760=', '96\u00d779='),\n    @('90\u00d745=', '13\u00d783='),\n    @('83\u00d773=', '72\u00d760='),\n    @('87\u00d724=', '50\u00d771='),\n    @('42\u00d795=', '17\u00d710='),\n    @('61\u00d736=', '61\u00d787='),\n    @('74\u00d745=', '43\u00d748='),\n    @('100\u00d730=', '80\u00d723='),\n    @('31\u00d726=', '22\u00d757='),\n    @('24\u00d770=', '18\u00d733='),\n    @('74\u00d751=', '80\u00d789='),\n    @('47\u00d744=', '13\u00d798='),\n    @('46\u00d727=', '10\u00d795='),\n    @('68\u00d727=', '22\u00d714='),\n    @('27\u00d770=', '91\u00d779='),\n    @('72\u00d777=', '85\u00d738='),\n    @('53\u00d732=', '95\u00d710='),\n    @('94\u00d770=', '11\u00d769='),\n    @('21\u00d748=', '41\u00d719='),\n    @('87\u00d743=', '79\u00d754='),\n    @('32\u00d783=', '86\u00d781='),\n    @('95\u00d714=', '94\u00d794='),\n    @('71\u00d738=', '69\u00d793='),\n    @('79\u00d798=', '30\u00d789='),\n    @('69\u00d741=', '76\u00d779='),\n    @('24\u00d740=', '74\u00d765='),\n    @('30\u00d761=', '61\u00d787='),\n    @('30\u00d773=', '13\u00d777='),\n    @('13\u00d769=', '24\u00d721='),\n    @('15\u00d721=', '73\u00d757='),\n    @('32\u00d768=', '74\u00d730='),\n    @('90\u00d762=', '24\u00d735='),\n    @('62\u00d747=', '65\u00d745='),\n    @('54\u00d726=', '30\u00d796='),\n    @('42\u00d750=', '72\u00d780='),\n    @('42\u00d729=', '63\u00d740='),\n    @('12\u00d734=', '28\u00d762='),\n    @('22\u00d734=', '38\u00d799='),\n    @('41\u00d760=', '93\u00d718='),\n    @('92\u00d782=', '81\u00d738='),\n    @('76\u00d776=', '47\u00d726='),\n    @('70\u00d732=', '19\u00d770='),\n    @('16\u00d753=', '59\u00d787='),\n    @('98\u00d757=', '95\u00d790='),\n    @('74\u00d738=', '42\u00d718='),\n    @('56\u00d725=', '30\u00d731='),\n    @('54\u00d7100=', '23\u00d765='),\n    @('52\u00d727=', '55\u00d764='),\n    @('44\u00d722=', '72\u00d786='),\n    @('87\u00d733=', '81\u00d726='),\n    @('11\u00d713=', '59\u00d758='),\n    @('89\u00d778=', '81\u00d739='),\n    @('38\u00d782=', '53\u00d748='),\n    @('50\u00d760=', '66\u00d756='),\n    @('42\u00d783=', '23\u00d763='),\n    @('11\u00d745=', '97\u00d787='),\n    @('67\u00d769=', '24\u00d743='),\n    @('12\u00d757=', '66\u00d721='),\n    @('22\u00d772=', '40\u00d790='),\n    @('44\u00d792=', '21\u00d733='),\n    @('11\u00d743=', '63\u00d751='),\n    @('55\u00d722=', '52\u00d767='),\n    @('83\u00d768=', '21\u00d736='),\n    @('30\u00d712=', '24\u00d714=')\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $p[0]\n    $find.Replacement.Text = $p[1]\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n}\n"}
